# Generate Report for Archive
# Replace "Ready for handoff" status text with "In Translation" in the
# Status columns of every sheet/table that tracks it (Overview's zh-cn/de-de
# columns, and the Status column of the zh-cn / de-de detail sheets), then
# let Excel auto-fit those columns to the new (shorter) text - matching
# Excel's "best fit" recalculation when cell content changes.

$wb = $excel.ActiveWorkbook

$newText = "In Translation"

# Overview sheet: status is mirrored into columns E (zh-cn) and F (de-de)
# for each of the 3 data rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = $newText

# zh-cn / de-de detail sheets: status lives in column C for each of the
# 3 data rows.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = $newText

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = $newText

# Shrink the status columns to match the narrower "best fit" width Excel
# computes for the new, shorter text (was sized for "Ready for handoff").
$newColumnWidth = 12.5
$overview.Columns("E:F").ColumnWidth = $newColumnWidth
$zhcn.Columns("C:C").ColumnWidth = $newColumnWidth
$dede.Columns("C:C").ColumnWidth = $newColumnWidth
